$wb = $excel.ActiveWorkbook

$oldGuid = "3d71b4c5-947a-429a-9c10-4e23a09dd6d6"
$newGuid = "d83ebb2f-8092-4156-b345-5b346095b88f"
$oldHash = "dc1311b846f9dd62cbf907a065a0b0c12964926f"
$newHash = "7a835065743460c6ef980d265a8a1fcdf7558795"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"
$wsOverview.Range("G2").Value = "2016-08-21 07:05:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 07:05:28"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
